$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet 1: LP1912  (columns: A blank, B Hora_Scrap, C Hora_Llegada,
#                   D Linea, E Minutos, F Parada, G Fecha)
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 17:32:38"
$ws1.Cells.Item(3,1).Value = "Total filas: 476"

$rows1 = @(
    @{ R=455; B="17:32:27"; C="17:35"; D="23_HERNANDEZ"; E=3; F="LP1912"; G="30/12/2025" }
    @{ R=456; B="17:32:27"; C="17:37"; D="27_EL RETIRO"; E=5; F="LP1912"; G="30/12/2025" }
    @{ R=457; B="17:32:27"; C="17:38"; D="17_ROMERO"; E=6; F="LP1912"; G="30/12/2025" }
    @{ R=458; B="17:32:27"; C="17:40"; D="16_SANTA ANA"; E=8; F="LP1912"; G="30/12/2025" }
    @{ R=459; B="17:32:27"; C="17:45"; D="15_ABASTO"; E=13; F="LP1912"; G="30/12/2025" }
    @{ R=460; B="17:32:27"; C="17:46"; D="10_OLMOS"; E=14; F="LP1912"; G="30/12/2025" }
    @{ R=461; B="17:32:27"; C="17:52"; D="81_EL PELIGRO"; E=20; F="LP1912"; G="30/12/2025" }
    @{ R=462; B="17:32:27"; C="18:00"; D="16_SANTA ANA"; E=28; F="LP1912"; G="30/12/2025" }
    @{ R=463; B="17:32:27"; C="18:04"; D="17_ROMERO"; E=32; F="LP1912"; G="30/12/2025" }
    @{ R=464; B="17:32:27"; C="18:04"; D="23_HERNANDEZ"; E=32; F="LP1912"; G="30/12/2025" }
    @{ R=465; B="17:32:27"; C="18:10"; D="16_SANTA ANA"; E=38; F="LP1912"; G="30/12/2025" }
    @{ R=466; B="17:32:27"; C="18:16"; D="10_OLMOS"; E=44; F="LP1912"; G="30/12/2025" }
    @{ R=467; B="17:32:27"; C="18:16"; D="15_ABASTO"; E=44; F="LP1912"; G="30/12/2025" }
    @{ R=468; B="17:32:27"; C="18:21"; D="26_HERNANDEZ"; E=49; F="LP1912"; G="30/12/2025" }
    @{ R=469; B="17:32:27"; C="18:24"; D="14_ABASTO"; E=52; F="LP1912"; G="30/12/2025" }
    @{ R=470; B="17:32:27"; C="18:28"; D="215C_EL PATO"; E=56; F="LP1912"; G="30/12/2025" }
    @{ R=471; B="17:32:27"; C="18:32"; D="11X44_ETCHEVERRY"; E=60; F="LP1912"; G="30/12/2025" }
    @{ R=472; B="17:32:27"; C="18:34"; D="23_HERNANDEZ"; E=62; F="LP1912"; G="30/12/2025" }
    @{ R=473; B="17:32:27"; C="18:40"; D="15_ABASTO"; E=68; F="LP1912"; G="30/12/2025" }
    @{ R=474; B="17:32:27"; C="18:48"; D="14X44_ABASTO"; E=76; F="LP1912"; G="30/12/2025" }
    @{ R=475; B="17:32:27"; C="18:56"; D="10_OLMOS"; E=84; F="LP1912"; G="30/12/2025" }
    @{ R=476; B="17:32:27"; C="18:59"; D="14_ABASTO"; E=87; F="LP1912"; G="30/12/2025" }
    @{ R=477; B="17:32:27"; C="19:04"; D="11_ETCHEVERRY"; E=92; F="LP1912"; G="30/12/2025" }
)

foreach ($row in $rows1) {
    $ws1.Cells.Item($row.R, 2).Value = $row.B
    $ws1.Cells.Item($row.R, 3).Value = $row.C
    $ws1.Cells.Item($row.R, 4).Value = $row.D
    $ws1.Cells.Item($row.R, 5).Value = $row.E
    $ws1.Cells.Item($row.R, 6).Value = $row.F
    $ws1.Cells.Item($row.R, 7).Value = $row.G
}

# -----------------------------------------------------------------
# Sheet 2: LP1912-215  (columns: A blank, B Fecha, C Hora_Scrap,
#                       D Hora_Llegada, E Linea, F Minutos, G Parada)
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 17:32:38"
$ws2.Cells.Item(3,1).Value = "Total filas: 31"

$rows2 = @(
    @{ R=32; B="30/12/2025"; C="17:32:27"; D="18:28"; E="215C_EL PATO"; F=56; G="LP1912" }
)

foreach ($row in $rows2) {
    $ws2.Cells.Item($row.R, 2).Value = $row.B
    $ws2.Cells.Item($row.R, 3).Value = $row.C
    $ws2.Cells.Item($row.R, 4).Value = $row.D
    $ws2.Cells.Item($row.R, 5).Value = $row.E
    $ws2.Cells.Item($row.R, 6).Value = $row.F
    $ws2.Cells.Item($row.R, 7).Value = $row.G
}

# -----------------------------------------------------------------
# Sheet 3: 6203-6173  (columns: A blank, B Fecha, C Hora_Scrap,
#                      D Hora_Llegada, E Linea, F Minutos, G Parada)
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 17:32:38"
$ws3.Cells.Item(3,1).Value = "Total filas: 62"

$rows3 = @(
    @{ R=62; B="30/12/2025"; C="17:32:32"; D="18:04"; E="215C_LA PLATA"; F=32; G="L6203" }
    @{ R=63; B="30/12/2025"; C="17:32:38"; D="18:52"; E="215A_LA PLATA"; F=80; G="L6173" }
)

foreach ($row in $rows3) {
    $ws3.Cells.Item($row.R, 2).Value = $row.B
    $ws3.Cells.Item($row.R, 3).Value = $row.C
    $ws3.Cells.Item($row.R, 4).Value = $row.D
    $ws3.Cells.Item($row.R, 5).Value = $row.E
    $ws3.Cells.Item($row.R, 6).Value = $row.F
    $ws3.Cells.Item($row.R, 7).Value = $row.G
}

Write-Host "Sheet1 rows:" $ws1.UsedRange.Rows.Count
Write-Host "Sheet2 rows:" $ws2.UsedRange.Rows.Count
Write-Host "Sheet3 rows:" $ws3.UsedRange.Rows.Count

